$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.655.67'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').Value = '2.287.72'
$ws.Range('E3').Value = '  +0.12%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '506.26'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '128.94'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.37%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.530'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('D9').Value = '2.310.90'
$ws.Range('E9').Value = '  +0.69%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0971'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.70%  '
$ws.Range('E11').Value = '  +1.72%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.342'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.61%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.93'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.46%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.59'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.49%  '
$ws.Range('D15').Value = '2.707.12'
$ws.Range('E15').Value = '  +0.60%  '
$ws.Range('D16').Value = '54.769.22'
$ws.Range('E16').Value = '  +0.81%  '
$ws.Range('E17').Value = '  +1.26%  '
$ws.Range('D18').Value = '2.282.51'
$ws.Range('E18').Value = '  -0.62%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.61'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.28%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.18'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.42%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.62'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.41%  '
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '308.28'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.97%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.36'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.78%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.996'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.26%  '
$ws.Range('E26').Value = '  -0.48%  '
$ws.Range('E27').Value = '  +2.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '171.48'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.49%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.12'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.94%  '
$ws.Range('D30').Value = '0.0₃0706'
$ws.Range('E30').Value = '  +2.13%  '
$ws.Range('E31').Value = '  +0.65%  '
$ws.Range('E32').Value = '  +5.19%  '
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.02'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.37%  '
$ws.Range('E35').Value = '  -0.29%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.906'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.29%  '
$ws.Range('E37').Value = '  +0.61%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.83'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.71%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.63'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.54%  '
$ws.Range('E40').Value = '  +0.53%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.43'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.81%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.08'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '131.94'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.59%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.42'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.65%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '252.60'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.89%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0500'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.75%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0912'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.88%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.554'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.80%  '
$ws.Range('E49').Value = '  +0.67%  '
$ws.Range('E50').Value = '  +0.63%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '10.82'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.41%  '
